# Scheduled-runner update: refresh Market Board price / profit columns (H:N)
# on the Leve-profit tables across sheets ALC, BSM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 476.7647
$ws.Range("I6").Value = 125.416664
$ws.Range("J6").Value = 1320
$ws.Range("K6").Value = 376.249992
$ws.Range("L6").Value = 3960
$ws.Range("M6").Value = -264.249992
$ws.Range("N6").Value = -4184
$ws.Range("H19").Value = 1021.5
$ws.Range("I19").Value = 795.0833
$ws.Range("K19").Value = 795.0833
$ws.Range("M19").Value = -620.0833
$ws.Range("H93").Value = 29750
$ws.Range("J93").Value = 29750
$ws.Range("L93").Value = 29750
$ws.Range("N93").Value = -34742
$ws.Range("H113").Value = 4249.3613
$ws.Range("I113").Value = 4361.25
$ws.Range("J113").Value = 4025.5833
$ws.Range("K113").Value = 4361.25
$ws.Range("L113").Value = 4025.5833
$ws.Range("M113").Value = -1107.25
$ws.Range("N113").Value = -10533.5833
$ws.Range("H116").Value = 2791.25
$ws.Range("I116").Value = 2721.6667
$ws.Range("K116").Value = 2721.6667
$ws.Range("M116").Value = 720.3332999999998
$ws.Range("H132").Value = 4794.4517
$ws.Range("I132").Value = 1593.6296
$ws.Range("J132").Value = 26400
$ws.Range("K132").Value = 4780.8888
$ws.Range("L132").Value = 79200
$ws.Range("M132").Value = -2250.8888
$ws.Range("N132").Value = -84260
$ws.Range("H138").Value = 1970
$ws.Range("I138").Value = 1260.4082
$ws.Range("J138").Value = 2885
$ws.Range("K138").Value = 3781.2246
$ws.Range("L138").Value = 8655
$ws.Range("M138").Value = 1358.7754
$ws.Range("N138").Value = -18935

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 27028844
$ws.Range("I20").Value = 41668252
$ws.Range("J20").Value = 2244.2307
$ws.Range("K20").Value = 41668252
$ws.Range("L20").Value = 2244.2307
$ws.Range("M20").Value = -41668005
$ws.Range("N20").Value = -2738.2307
$ws.Range("H92").Value = 21950.25
$ws.Range("J92").Value = 21950.25
$ws.Range("L92").Value = 21950.25
$ws.Range("N92").Value = -26942.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 19466.6
$ws.Range("I86").Value = 13083.25
$ws.Range("J86").Value = 45000
$ws.Range("K86").Value = 13083.25
$ws.Range("L86").Value = 45000
$ws.Range("M86").Value = -11960.25
$ws.Range("N86").Value = -47246
$ws.Range("H89").Value = 19466.6
$ws.Range("I89").Value = 13083.25
$ws.Range("J89").Value = 45000
$ws.Range("K89").Value = 65416.25
$ws.Range("L89").Value = 225000
$ws.Range("M89").Value = -59800.25
$ws.Range("N89").Value = -236232
$ws.Range("H125").Value = 19818.182
$ws.Range("J125").Value = 19818.182
$ws.Range("L125").Value = 19818.182
$ws.Range("N125").Value = -24738.182
$ws.Range("H132").Value = 1196867
$ws.Range("I132").Value = 1812.7368
$ws.Range("J132").Value = 3089036.2
$ws.Range("K132").Value = 5438.2104
$ws.Range("L132").Value = 9267108.600000001
$ws.Range("M132").Value = -2908.2104
$ws.Range("N132").Value = -9272168.600000001
$ws.Range("H138").Value = 36577.9
$ws.Range("J138").Value = 36577.9
$ws.Range("L138").Value = 36577.9
$ws.Range("N138").Value = -46857.9
$ws.Range("H141").Value = 87468.25
$ws.Range("J141").Value = 87468.25
$ws.Range("L141").Value = 87468.25
$ws.Range("N141").Value = -97828.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 7658
$ws.Range("I120").Value = 6572.5
$ws.Range("K120").Value = 19717.5
$ws.Range("M120").Value = -14879.5
$ws.Range("H129").Value = 1798.2162
$ws.Range("I129").Value = 1019.41174
$ws.Range("J129").Value = 2460.2
$ws.Range("K129").Value = 3058.23522
$ws.Range("L129").Value = 7380.599999999999
$ws.Range("M129").Value = 1941.76478
$ws.Range("N129").Value = -17380.6
$ws.Range("H131").Value = 935.0599999999999
$ws.Range("I131").Value = 269.5
$ws.Range("J131").Value = 977.54254
$ws.Range("K131").Value = 808.5
$ws.Range("L131").Value = 2932.62762
$ws.Range("M131").Value = 4231.5
$ws.Range("N131").Value = -13012.62762
$ws.Range("H132").Value = 2158.0952
$ws.Range("I132").Value = 2541.3333
$ws.Range("K132").Value = 22871.9997
$ws.Range("M132").Value = -20341.9997
$ws.Range("H140").Value = 2372.0293
$ws.Range("I140").Value = 929.0833
$ws.Range("J140").Value = 3159.0908
$ws.Range("K140").Value = 2787.2499
$ws.Range("L140").Value = 9477.2724
$ws.Range("M140").Value = 2392.7501
$ws.Range("N140").Value = -19837.2724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5464.96
$ws.Range("I70").Value = 5102.125
$ws.Range("J70").Value = 6110
$ws.Range("K70").Value = 5102.125
$ws.Range("L70").Value = 6110
$ws.Range("M70").Value = -4832.125
$ws.Range("N70").Value = -6650
$ws.Range("H73").Value = 5464.96
$ws.Range("I73").Value = 5102.125
$ws.Range("J73").Value = 6110
$ws.Range("K73").Value = 5102.125
$ws.Range("L73").Value = 6110
$ws.Range("M73").Value = -4166.125
$ws.Range("N73").Value = -7982
$ws.Range("H80").Value = 2229.8518
$ws.Range("I80").Value = 2244.4546
$ws.Range("J80").Value = 2219.8125
$ws.Range("K80").Value = 2244.4546
$ws.Range("L80").Value = 2219.8125
$ws.Range("M80").Value = -1246.4546
$ws.Range("N80").Value = -4215.8125
$ws.Range("H83").Value = 2229.8518
$ws.Range("I83").Value = 2244.4546
$ws.Range("J83").Value = 2219.8125
$ws.Range("K83").Value = 11222.273
$ws.Range("L83").Value = 11099.0625
$ws.Range("M83").Value = -6230.273000000001
$ws.Range("N83").Value = -21083.0625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2941.2778
$ws.Range("I40").Value = 3111.5833
$ws.Range("J40").Value = 2600.6667
$ws.Range("K40").Value = 3111.5833
$ws.Range("L40").Value = 2600.6667
$ws.Range("M40").Value = -2975.5833
$ws.Range("N40").Value = -2872.6667
$ws.Range("H82").Value = 1103.0769
$ws.Range("I82").Value = 1845
$ws.Range("J82").Value = 968.1818
$ws.Range("K82").Value = 1845
$ws.Range("L82").Value = 968.1818
$ws.Range("M82").Value = -1484
$ws.Range("N82").Value = -1690.1818
$ws.Range("H85").Value = 1103.0769
$ws.Range("I85").Value = 1845
$ws.Range("J85").Value = 968.1818
$ws.Range("K85").Value = 1845
$ws.Range("L85").Value = 968.1818
$ws.Range("M85").Value = -597
$ws.Range("N85").Value = -3464.1818
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H122").Value = 2740.4167
$ws.Range("I122").Value = 1956
$ws.Range("K122").Value = 5868
$ws.Range("M122").Value = -3418

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 28929.8
$ws.Range("J92").Value = 28929.8
$ws.Range("L92").Value = 28929.8
$ws.Range("N92").Value = -33921.8
